$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Objetivos: (row 10) gains its body text in B and C (was a stray docente name)
$t10 = "1. Introduzir e discutir conceitos e técnicas estatísticas para controle e melhoria da qualidade de produtos fabricados e processos de fabricação;2 Fornecer subsídios para que o aluno tenha condições de utilizar essas técnicas e conceitos na sua vida profissional futura."
$ws.Range("B10").Value = $t10
$ws.Range("C10").Value = $t10

# 2) Insert a blank row at 13 - pushes old rows 13-23 down to 14-24,
#    carrying their row heights and formatting with them.
$ws.Rows.Item(13).Insert()

# The inserted row 13 picked up a stray empty A13 cell (bold style) from the
# row above via autofill; the target row 13 has no A cell at all, so drop it.
$ws.Range("A13").Clear()

# 3) Fill the newly-revealed "label-only" rows with their B/C content
$t13 = "5840917 - Fabricio Maciel Gomes"
$ws.Range("B13").Value = $t13
$ws.Range("C13").Value = $t13
# match B13/C13 formatting (wrap-text / red-wrap-text) to the rest of the table
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$t14 = "Sistemas de Medição, Fundamentos do Controle Estatístico da Qualidade e do Processo, Gráficos de Controle por Variáveis, Gráficos de Controle por Atributos, Gráficos de Controle para Processos Auto-correlacionados, Analise de Capacidade do Processo, Inspeção da Qualidade, Estudos de casos."
$ws.Range("B14").Value = $t14
$ws.Range("C14").Value = $t14

$t16 = "1. Sistemas de Medição.1.1. Planejamento do Sistema de Medição;1.2. Impacto da Variabilidade do Sistema de Medição no Produto;1.3. Sistemas de Medição por Atributos;1.4. Tendência e Linearidade;1.5. Análise de Repetitividade e Reprodutibilidade;2. Fundamentos do Controle Estatístico da Qualidade e do Processo.2.1. Importância do Controle Estatístico da Qualidade e do Processo;2.2. Naturezas das Variações;2.3. Causas Comuns e Causas Especiais de Variações;3. Gráficos de Controle por Variáveis3.1. Gráficos de Controle por Médias;3.2. Gráficos de Controle por Amplitude;3.3. Gráficos de Controle por Desvio Padrão;3.4. Análise de Desempenho dos Gráficos de Controle por Variáveis;4. Gráficos de Controle por Atributos4.1. Gráficos de Controle por Número de Não Conformidades;4.2. Gráficos de Controle por Fração Não Conforme;4.3. Gráficos de Controle por Número de Defeitos4.4. Gráficos de Controle por Não Conformidades por Amostra;5. Gráficos de Controle para Processos Auto-correlacionados5.1. Gráficos de Controle por Amplitude Móvel;5.2. Gráficos de Controle por Soma Acumulada (CUSUM).5.3. Gráficos de Controle por Média Móvel Ponderada Exponencialmente (EWMA)6. Analise de Capacidade do Processo6.1. Índices de Capacidade do Processo;6.2. Índices de Performance do Processo; 7. Inspeção da Qualidade7.1. Planos de Amostragem7.2. Inspeção para Aceitação;7.3. Inspeção Retificadora;8. Estudos de casos"
$ws.Range("B16").Value = $t16
$ws.Range("C16").Value = $t16

$t19 = "Aulas expositivas teóricas, aulas práticas, aulas de laboratório, aulas de exercícios."
$ws.Range("B19").Value = $t19
$ws.Range("C19").Value = $t19

$t20 = "MF = (0,5*P1 + 0,5*P2), onde P1 e P2 são provas."
$ws.Range("B20").Value = $t20
$ws.Range("C20").Value = $t20

$t21 = "Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação."
$ws.Range("B21").Value = $t21
$ws.Range("C21").Value = $t21

$t22 = "1. COSTA, A.F.B., EPPRECHT, E.K., CARPINETTI, L.C.R., Controle Estatístico da Qualidade, 2ª ed., Editora Atlas, 2005.`n2. MONTGOMERY, D.C., Introdução ao Controle Estatístico da Qualidade, 4ª ed., Livros Técnicos e Científicos, 2004.`n3. GRANT, E., LEAVENWORTH, R., Statistical Quality Control, 7ªed., McGraw-Hill, 1996.`n4. WERKENA, M.C.C., Ferramentas Estatísticas Básicas para o Gerenciamento de Processos, Editora FCO, 1996."
$ws.Range("B22").Value = $t22
$ws.Range("C22").Value = $t22

